# Apply the tracked edit from the commit:
#   slide 10 (sldId 762) - shape id 5 "Smiley Face 4": set the Alt-Text
#   "Description" field to "QuizAnswer". In OOXML this sets the
#   p:cNvPr/@descr attribute, i.e. <p:cNvPr id="5" name="Smiley Face 4"
#   descr="QuizAnswer"/>.
#
# (The footer "today" date fields that PowerPoint auto-recached on the
# slide master / layouts / notes master during this same save are an
# incidental side effect of autosave on a later day, not an authored
# content edit, and PowerPoint's object model does not expose a way to
# rewrite a live datetimeFigureOut field's cached text without
# collapsing it to a plain (non-updating) run, so they are intentionally
# left untouched here.)

$p = $ppt.ActivePresentation

$s = $p.Slides.Item(10)

$shp = $s.Shapes.Item(4)

if ($shp.Name -ne "Smiley Face 4") {
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        if ($s.Shapes.Item($i).Name -eq "Smiley Face 4") {
            $shp = $s.Shapes.Item($i)
        }
    }
}

$shp.AlternativeText = "QuizAnswer"
